$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-6), including a new row 6
$data = @(
    @(1, 7, 9, 2, 8, -5, -1, 12, 5),
    @(2, 5, 9, 1, 7, -4, -2, 23, 5),
    @(3, 9, 5, 8, 0, -1, -5, 56, 5),
    @(4, 8, 8, 5, 5, -3, -3, 34, 5),
    @(5, 6, 6, 4, 2, -2, -4, 45, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$ws.Range("I1").Select()
